$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column Q (rows 30-33): new "interessante links" hyperlinks ---
# (added first so the shared-string table fills in the same order as the
# original authoring session: links first, then the new log entry text)
$ws.Cells.Item(30, 17).Value = "https://stackoverflow.com/questions/48492993/firestore-get-documentsnapshots-fields-value"
$ws.Hyperlinks.Add($ws.Cells.Item(30, 17), "https://stackoverflow.com/questions/48492993/firestore-get-documentsnapshots-fields-value")
$ws.Cells.Item(30, 17).Style = "Hyperlink"

$ws.Cells.Item(31, 17).Value = "https://stackoverflow.com/questions/35805891/how-to-get-only-even-numbers-from-list"
$ws.Hyperlinks.Add($ws.Cells.Item(31, 17), "https://stackoverflow.com/questions/35805891/how-to-get-only-even-numbers-from-list")
$ws.Cells.Item(31, 17).Style = "Hyperlink"

$ws.Cells.Item(32, 17).Value = "https://www.journaldev.com/12478/android-searchview-example-tutorial"
$ws.Hyperlinks.Add($ws.Cells.Item(32, 17), "https://www.journaldev.com/12478/android-searchview-example-tutorial")
$ws.Cells.Item(32, 17).Style = "Hyperlink"

$ws.Cells.Item(33, 17).Value = "https://abhiandroid.com/ui/searchview"
$ws.Hyperlinks.Add($ws.Cells.Item(33, 17), "https://abhiandroid.com/ui/searchview")
$ws.Cells.Item(33, 17).Style = "Hyperlink"

# --- Row 30: new logboek entry (20/04/2019, "1 uur 50 minuten", description) ---
$ws.Cells.Item(30, 1).Value = 43575
$ws.Cells.Item(30, 1).NumberFormat = "d-mmm"

$ws.Cells.Item(30, 2).Value = "1 uur 50 minuten"

$ws.Cells.Item(30, 3).Value = "route beschrijvingen ophalen uit firestore en in cardviews steken, uitzoeken hoe een searchview werkt"

# --- Selection cursor left where the author's last edit landed ---
$ws.Range("D32").Select()
